$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows right before the current row 633, shifting rows 633:713
# down to 640:720 (Excel's default Insert shifts cells down and carries the
# formatting of the row above into the freshly inserted rows).
$ws.Range("A633:T639").Insert()

# Fill in the 7 newly inserted rows (633-639) with their data.
$newRows = @(
    @{ Row = 633; D = 44474; K = "Edranol";           L = "Primera";             M = 50;  N = 3000; O = 3000; P = 3000; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Región de O'Higgins";  S = 3000; T = 1 },
    @{ Row = 634; D = 44474; K = "Hass";               L = "1a nueva(o)";         M = 100; N = 3500; O = 3500; P = 3500; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 3500; T = 1 },
    @{ Row = 635; D = 44474; K = "Hass";               L = "2a nueva(o)";         M = 120; N = 2800; O = 2800; P = 2800; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 2800; T = 1 },
    @{ Row = 636; D = 44474; K = "Hass";               L = "4a nueva (o)";        M = 50;  N = 1500; O = 1500; P = 1500; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 1500; T = 1 },
    @{ Row = 637; D = 44474; K = "Hass";               L = "Especial nueva (o)";  M = 50;  N = 3600; O = 3600; P = 3600; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 3600; T = 1 },
    @{ Row = 638; D = 44474; K = "Hass";               L = "Segunda";             M = 200; N = 3000; O = 3000; P = 3000; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 3000; T = 1 },
    @{ Row = 639; D = 44474; K = "Negra de La Cruz";   L = "Primera";             M = 50;  N = 3000; O = 3000; P = 3000; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Región de O'Higgins";  S = 3000; T = 1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value = "La Araucanía"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100106
    $ws.Cells.Item($row, 8).Value = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
